$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 07:05"

# Row 11 - Alemania: update active cases / recovered
$ws.Range("D11").Value = 154600
$ws.Range("E11").Value = 14002

# Row 73 - Tailandia: update totals/new/active/recovered
$ws.Range("B73").Value = 3031
$ws.Range("C73").Value = 3
$ws.Range("D73").Value = 2857
$ws.Range("E73").Value = 118

# Rows 78-80 - Sudan's case counts grew enough to move above Senegal and
# Bosnia y Herzegovina in the (descending, by total cases) ranking, so the
# three countries shift position while their underlying data updates.
# Row 78 becomes Sudan (previously Senegal) with refreshed figures.
$ws.Range("A78").Value = "Sudan"
$ws.Range("B78").Value = 2591
$ws.Range("C78").Value = 302
$ws.Range("D78").Value = 247
$ws.Range("E78").Value = 2239
$ws.Range("G78").Value = 8
$ws.Range("H78").Value = 105

# Row 79 becomes Senegal (previously Bosnia y Herzegovina), carrying the
# figures that used to belong to Senegal in row 78.
$ws.Range("A79").Value = "Senegal"
$ws.Range("B79").Value = 2480
$ws.Range("D79").Value = 973
$ws.Range("E79").Value = 1482
$ws.Range("H79").Value = 25

# Row 80 becomes Bosnia y Herzegovina (previously Sudan), carrying the
# figures that used to belong to Bosnia y Herzegovina in row 79.
$ws.Range("A80").Value = "Bosnia y Herzegovina"
$ws.Range("B80").Value = 2290
$ws.Range("D80").Value = 1436
$ws.Range("E80").Value = 721
$ws.Range("H80").Value = 133
